$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password/manager values shared across rows 2-4
$ws.Range("B2:B4").Value = "mngr126739"
$ws.Range("C2:C4").Value = "netyveb"

# Update the selected cell
$ws.Range("D7").Select()
